# Regenerate merged AHB files
# - Rename the header row labels from *_old / *_new to *_FV2210 / *_FV2304
# - Turn the data range into an Excel Table (ListObject)
# - Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header cells (row 1, columns A:U)
$lastCol = 21
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $label = $cell.Text
    if ($label -like "*_old") {
        $cell.Value = ($label -replace "_old$", "_FV2210")
    } elseif ($label -like "*_new") {
        $cell.Value = ($label -replace "_new$", "_FV2304")
    }
}

# 2) Wrap the used range in a native Excel table, re-using the renamed headers
$dataRange = $ws.Range("A1:U80")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"

# 3) Freeze panes above row 2 (i.e. freeze the header row)
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
